{"js": "// Replace the 25 \"a\u00f7b=c, d\" answer strings in the worksheet table with the\n// new values below, in document order. Every populated table cell changes;\n// empty filler cells (used for spacing between problem rows) are skipped.\n// A plain content-based search/replace is unsafe here because one source\n// string (\"376\u00f73=125, 1\") occurs twice and maps to two different targets,\n// so we walk the table cells positionally instead.\n\nconst replacements = [\n  \"360\u00f75=72, 0\",\n  \"995\u00f78=124, 3\",\n  \"936\u00f76=156, 0\",\n  \"807\u00f79=89, 6\",\n  \"864\u00f76=144, 0\",\n  \"352\u00f76=58, 4\",\n  \"949\u00f74=237, 1\",\n  \"643\u00f76=107, 1\",\n  \"158\u00f72=79, 0\",\n  \"265\u00f74=66, 1\",\n  \"511\u00f76=85, 1\",\n  \"424\u00f76=70, 4\",\n  \"538\u00f78=67, 2\",\n  \"523\u00f75=104, 3\",\n  \"791\u00f72=395, 1\",\n  \"659\u00f76=109, 5\",\n  \"366\u00f77=52, 2\",\n  \"310\u00f74=77, 2\",\n  \"840\u00f77=120, 0\",\n  \"636\u00f76=106, 0\",\n  \"726\u00f72=363, 0\",\n  \"121\u00f73=40, 1\",\n  \"113\u00f73=37, 2\",\n  \"163\u00f74=40, 3\",\n  \"449\u00f73=149, 2\",\n];\n\nconst table = context.document.body.tables.getFirst();\ntable.load(\"rowCount,values\");\nawait context.sync();\n\n// Find the (row, col) of every populated cell, in row-major document\n// order (this matches the order the values above should be applied in).\nconst targets = [];\nfor (let r = 0; r < table.rowCount; r++) {\n  const rowVals = table.values[r];\n  for (let c = 0; c < rowVals.length; c++) {\n    if (rowVals[c] !== \"\") {\n      targets.push({ r, c });\n    }\n  }\n}\n\nif (targets.length !== replacements.length) {\n  throw new Error(\n    \"Expected \" + replacements.length + \" populated cells, found \" + targets.length\n  );\n}\n\nfor (let i = 0; i < targets.length; i++) {\n  const { r, c } = targets[i];\n  const cell = table.getCell(r, c);\n  cell.body.getRange().insertText(replacements[i], Word.InsertLocation.replace);\n}\n\nawait context.sync();\n", "ps1": "# Replace the 25 \"a\u00f7b=c, d\" answer strings in the worksheet table with the\n# new values below, in document order. Every populated table cell changes;\n# empty filler cells (used for spacing between problem rows) are skipped.\n# A plain Find/Replace pass is unsafe here because one source string\n# (\"376\u00f73=125, 1\") occurs twice and maps to two different targets, so we\n# walk the table cells positionally instead.\n\n$replacements = @(\n  \"360\u00f75=72, 0\",\n  \"995\u00f78=124, 3\",\n  \"936\u00f76=156, 0\",\n  \"807\u00f79=89, 6\",\n  \"864\u00f76=144, 0\",\n  \"352\u00f76=58, 4\",\n  \"949\u00f74=237, 1\",\n  \"643\u00f76=107, 1\",\n  \"158\u00f72=79, 0\",\n  \"265\u00f74=66, 1\",\n  \"511\u00f76=85, 1\",\n  \"424\u00f76=70, 4\",\n  \"538\u00f78=67, 2\",\n  \"523\u00f75=104, 3\",\n  \"791\u00f72=395, 1\",\n  \"659\u00f76=109, 5\",\n  \"366\u00f77=52, 2\",\n  \"310\u00f74=77, 2\",\n  \"840\u00f77=120, 0\",\n  \"636\u00f76=106, 0\",\n  \"726\u00f72=363, 0\",\n  \"121\u00f73=40, 1\",\n  \"113\u00f73=37, 2\",\n  \"163\u00f74=40, 3\",\n  \"449\u00f73=149, 2\"\n)\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n$i = 0\nfor ($r = 1; $r -le $t.Rows.Count; $r++) {\n  for ($c = 1; $c -le $t.Columns.Count; $c++) {\n    $cell = $t.Cell($r, $c)\n    $cellRange = $cell.Range\n    # Drop the trailing end-of-cell marker(s) so we only look at real text.\n    $txt = $cellRange.Text.TrimEnd([char]13, [char]7)\n    if ($txt.Length -gt 0) {\n      if ($i -ge $replacements.Length) {\n        throw \"More populated cells than replacement values\"\n      }\n      $cellRange.Text = $replacements[$i]\n      $i++\n    }\n  }\n}\n\nif ($i -ne $replacements.Length) {\n  throw \"Expected $($replacements.Length) populated cells, found $i\"\n}\n"}
